$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 54) right below the existing data, which
# currently ends at row 53.
$newRow = 54
$srcRow = $newRow - 1

# Write the new values. Setting the date as plain text first (without first
# forcing a Text number format) would make Excel auto-convert "2026/01/03"
# into a date serial, so the number format is matched up afterwards by
# copying the previous row's formatting.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026/01/03"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1130

# Copy the formatting (cell style / alignment) from the previous row so the
# new row matches the rest of the table's look.
$ws.Range("A$srcRow`:C$srcRow").Copy()
$ws.Range("A$newRow`:C$newRow").PasteSpecial(-4122)

$wb.Save()
